# Epoch Accuracy.xlsx — refresh the per-epoch accuracy numbers (column B)
# with the results from the latest training run ("FO3 Froze Encoder 12345"),
# and refresh the stale Python object-repr strings in column A that carry
# the new run's memory address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated per-epoch accuracy values -------------------------------------
$ws.Range("B5").Value   = 0.984375
$ws.Range("B6").Value   = 0.96875
$ws.Range("B7").Value   = 1
$ws.Range("B8").Value   = 0.953125
$ws.Range("B9").Value   = 0.90625
$ws.Range("B10").Value  = 0.921875
$ws.Range("B12").Value  = 0.890625
$ws.Range("B13").Value  = 0.875
$ws.Range("B14").Value  = 0.859375
$ws.Range("B15").Value  = 0.8125
$ws.Range("B17").Value  = 0.828125
$ws.Range("B18").Value  = 0.859375
$ws.Range("B19").Value  = 0.875
$ws.Range("B20").Value  = 0.78125
$ws.Range("B21").Value  = 0.796875
$ws.Range("B22").Value  = 0.765625
$ws.Range("B23").Value  = 0.796875
$ws.Range("B24").Value  = 0.765625
$ws.Range("B25").Value  = 0.78125
$ws.Range("B26").Value  = 0.796875
$ws.Range("B28").Value  = 0.8125
$ws.Range("B30").Value  = 0.828125
$ws.Range("B31").Value  = 0.828125
$ws.Range("B32").Value  = 0.828125
$ws.Range("B33").Value  = 0.828125
$ws.Range("B34").Value  = 0.828125
$ws.Range("B35").Value  = 0.828125
$ws.Range("B36").Value  = 0.828125
$ws.Range("B37").Value  = 0.8125
$ws.Range("B46").Value  = 0.796875
$ws.Range("B47").Value  = 0.796875
$ws.Range("B48").Value  = 0.796875
$ws.Range("B49").Value  = 0.796875
$ws.Range("B50").Value  = 0.796875
$ws.Range("B54").Value  = 0.8125
$ws.Range("B59").Value  = 0.8125
$ws.Range("B65").Value  = 0.796875
$ws.Range("B66").Value  = 0.8125
$ws.Range("B67").Value  = 0.8125
$ws.Range("B68").Value  = 0.8125
$ws.Range("B69").Value  = 0.8125
$ws.Range("B70").Value  = 0.8125
$ws.Range("B73").Value  = 0.796875
$ws.Range("B74").Value  = 0.796875
$ws.Range("B75").Value  = 0.796875
$ws.Range("B76").Value  = 0.796875
$ws.Range("B77").Value  = 0.796875
$ws.Range("B78").Value  = 0.796875
$ws.Range("B79").Value  = 0.796875
$ws.Range("B80").Value  = 0.796875
$ws.Range("B81").Value  = 0.796875
$ws.Range("B82").Value  = 0.796875
$ws.Range("B83").Value  = 0.8125
$ws.Range("B84").Value  = 0.8125
$ws.Range("B85").Value  = 0.8125
$ws.Range("B86").Value  = 0.8125
$ws.Range("B87").Value  = 0.8125
$ws.Range("B88").Value  = 0.8125
$ws.Range("B89").Value  = 0.8125
$ws.Range("B90").Value  = 0.8125
$ws.Range("B91").Value  = 0.8125
$ws.Range("B92").Value  = 0.8125
$ws.Range("B93").Value  = 0.8125
$ws.Range("B94").Value  = 0.8125
$ws.Range("B95").Value  = 0.8125
$ws.Range("B96").Value  = 0.796875
$ws.Range("B97").Value  = 0.796875
$ws.Range("B98").Value  = 0.796875
$ws.Range("B99").Value  = 0.796875
$ws.Range("B100").Value = 0.796875
$ws.Range("B101").Value = 0.796875
$ws.Range("B102").Value = 0.796875
$ws.Range("B103").Value = 0.90625
$ws.Range("B104").Value = 0.796875
$ws.Range("B105").Value = 0.796875
$ws.Range("B106").Value = 0.796875
$ws.Range("B107").Value = 0.78125
$ws.Range("B108").Value = 0.8125
$ws.Range("B109").Value = 0.8125
$ws.Range("B110").Value = 0.78125
$ws.Range("B111").Value = 0.859375
$ws.Range("B112").Value = 0.875
$ws.Range("B113").Value = 0.703125
$ws.Range("B114").Value = 0.765625
$ws.Range("B115").Value = 0.84375
$ws.Range("B116").Value = 0.84375
$ws.Range("B117").Value = 0.859375
$ws.Range("B118").Value = 0.7213114754098361

# --- Refresh the DisplayOutputs repr string (new interpreter -> new id()) --
$newRepr = "<__main__.DisplayOutputs object at 0x7fccb2c1cca0>"
$ws.Range("A102:A118").Value = $newRepr

# --- Mirror the author's final on-screen state: select the whole sheet,
#     with the cursor left resting around D116 (bottom of the visible pane).
$ws.Range("D116").Select() | Out-Null
$ws.Cells.Select() | Out-Null
